# Applies the weekly Fruit/Vegetable price-table corrections for
# "Vega Modelo de Temuco - Papa": revised values for rows 651-663 and
# three newly-reported price records appended as rows 664-666.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct existing rows 651-663 ---

# Row 651
$ws.Range("D651").Value = 44656
$ws.Range("H651").Value = 'Patagonia'
$ws.Range("I651").Value = '1a (cosecha)'
$ws.Range("J651").Value = 300

# Row 652
$ws.Range("D652").Value = 44656
$ws.Range("H652").Value = 'Patagonia'
$ws.Range("I652").Value = '1a (cosecha)'
$ws.Range("J652").Value = 450

# Row 653
$ws.Range("D653").Value = 44656
$ws.Range("I653").Value = '1a (cosecha)'
$ws.Range("J653").Value = 380
$ws.Range("K653").Value = 7000
$ws.Range("L653").Value = 7000
$ws.Range("M653").Value = 7000
$ws.Range("N653").Value = '$/saco 25 kilos'
$ws.Range("P653").Value = 280

# Row 654
$ws.Range("H654").Value = 'Asterix'
$ws.Range("J654").Value = 800
$ws.Range("K654").Value = 8000
$ws.Range("L654").Value = 8000
$ws.Range("M654").Value = 8000
$ws.Range("N654").Value = '$/malla 25 kilos'
$ws.Range("P654").Value = 320

# Row 655
$ws.Range("D655").Value = 44432
$ws.Range("J655").Value = 200
$ws.Range("K655").Value = 7000
$ws.Range("L655").Value = 7000
$ws.Range("M655").Value = 7000
$ws.Range("N655").Value = '$/saco 25 kilos'
$ws.Range("P655").Value = 280

# Row 656
$ws.Range("D656").Value = 44432
$ws.Range("H656").Value = 'Rosara'
$ws.Range("I656").Value = '1a (guarda)'
$ws.Range("J656").Value = 600
$ws.Range("K656").Value = 8000
$ws.Range("L656").Value = 8000
$ws.Range("M656").Value = 8000
$ws.Range("N656").Value = '$/malla 25 kilos'
$ws.Range("P656").Value = 320

# Row 657
$ws.Range("D657").Value = 44432
$ws.Range("H657").Value = 'Rosara'
$ws.Range("J657").Value = 300
$ws.Range("K657").Value = 7000
$ws.Range("L657").Value = 7000
$ws.Range("M657").Value = 7000
$ws.Range("N657").Value = '$/saco 25 kilos'
$ws.Range("P657").Value = 280

# Row 658
$ws.Range("D658").Value = 44515
$ws.Range("H658").Value = 'Asterix'
$ws.Range("I658").Value = '1a (guarda)'
$ws.Range("J658").Value = 480
$ws.Range("K658").Value = 8000
$ws.Range("L658").Value = 8000
$ws.Range("M658").Value = 8000
$ws.Range("N658").Value = '$/malla 25 kilos'
$ws.Range("P658").Value = 320

# Row 659
$ws.Range("D659").Value = 44515
$ws.Range("H659").Value = 'Pehuenche'
$ws.Range("I659").Value = '1a nueva(o)'
$ws.Range("J659").Value = 850
$ws.Range("K659").Value = 13000
$ws.Range("L659").Value = 13000
$ws.Range("M659").Value = 13000
$ws.Range("N659").Value = '$/saco 25 kilos'
$ws.Range("P659").Value = 520

# Row 660
$ws.Range("D660").Value = 44508
$ws.Range("H660").Value = 'Asterix'
$ws.Range("I660").Value = '1a (guarda)'
$ws.Range("J660").Value = 900
$ws.Range("K660").Value = 8000
$ws.Range("L660").Value = 9000
$ws.Range("M660").Value = 8444
$ws.Range("P660").Value = 338

# Row 661
$ws.Range("D661").Value = 44508
$ws.Range("H661").Value = 'Pehuenche'
$ws.Range("J661").Value = 80
$ws.Range("K661").Value = 13000
$ws.Range("L661").Value = 13000
$ws.Range("M661").Value = 13000
$ws.Range("P661").Value = 520

# Row 662
$ws.Range("D662").Value = 44508
$ws.Range("H662").Value = 'Rosara'
$ws.Range("I662").Value = '1a (guarda)'
$ws.Range("J662").Value = 450
$ws.Range("K662").Value = 8000
$ws.Range("L662").Value = 9000
$ws.Range("M662").Value = 8444
$ws.Range("P662").Value = 338

# Row 663
$ws.Range("H663").Value = 'Patagonia'
$ws.Range("J663").Value = 500
$ws.Range("K663").Value = 7000
$ws.Range("L663").Value = 7000
$ws.Range("M663").Value = 7000
$ws.Range("N663").Value = '$/malla 25 kilos'
$ws.Range("P663").Value = 280

# --- Append new rows 664-666 ---

# Row 664
$ws.Range("D664").NumberFormat = $ws.Range("D663").NumberFormat
$ws.Range("A664").Value = 10
$ws.Range("B664").Value = 'Vega Modelo de Temuco'
$ws.Range("C664").Value = 'La Araucanía'
$ws.Range("D664").Value = 44592
$ws.Range("E664").Value = 9
$ws.Range("F664").Value = 100114001
$ws.Range("G664").Value = 'Papa'
$ws.Range("H664").Value = 'Patagonia'
$ws.Range("I664").Value = '1a nueva(o)'
$ws.Range("J664").Value = 500
$ws.Range("K664").Value = 6000
$ws.Range("L664").Value = 6000
$ws.Range("M664").Value = 6000
$ws.Range("N664").Value = '$/saco 25 kilos'
$ws.Range("O664").Value = 'Provincia de Cautín'
$ws.Range("P664").Value = 240
$ws.Range("Q664").Value = 25
$ws.Range("R664").Value = 'Hortaliza'

# Row 665
$ws.Range("D665").NumberFormat = $ws.Range("D664").NumberFormat
$ws.Range("A665").Value = 10
$ws.Range("B665").Value = 'Vega Modelo de Temuco'
$ws.Range("C665").Value = 'La Araucanía'
$ws.Range("D665").Value = 44592
$ws.Range("E665").Value = 9
$ws.Range("F665").Value = 100114001
$ws.Range("G665").Value = 'Papa'
$ws.Range("H665").Value = 'Rodeo'
$ws.Range("I665").Value = '1a nueva(o)'
$ws.Range("J665").Value = 600
$ws.Range("K665").Value = 7000
$ws.Range("L665").Value = 7000
$ws.Range("M665").Value = 7000
$ws.Range("N665").Value = '$/malla 25 kilos'
$ws.Range("O665").Value = 'Provincia de Cautín'
$ws.Range("P665").Value = 280
$ws.Range("Q665").Value = 25
$ws.Range("R665").Value = 'Hortaliza'

# Row 666
$ws.Range("D666").NumberFormat = $ws.Range("D665").NumberFormat
$ws.Range("A666").Value = 10
$ws.Range("B666").Value = 'Vega Modelo de Temuco'
$ws.Range("C666").Value = 'La Araucanía'
$ws.Range("D666").Value = 44592
$ws.Range("E666").Value = 9
$ws.Range("F666").Value = 100114001
$ws.Range("G666").Value = 'Papa'
$ws.Range("H666").Value = 'Rodeo'
$ws.Range("I666").Value = '1a nueva(o)'
$ws.Range("J666").Value = 900
$ws.Range("K666").Value = 6000
$ws.Range("L666").Value = 6000
$ws.Range("M666").Value = 6000
$ws.Range("N666").Value = '$/saco 25 kilos'
$ws.Range("O666").Value = 'Provincia de Cautín'
$ws.Range("P666").Value = 240
$ws.Range("Q666").Value = 25
$ws.Range("R666").Value = 'Hortaliza'
